$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Add a new "2022-Q4" sheet, placed right before the "2022-Q3" sheet.
#    We duplicate the existing "2022-Q3" sheet (same column layout /
#    header / styling) and then overwrite the numbers that changed.
# ------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)                      # new copy is inserted immediately before $q3,
                                    # named "2022-Q3 (2)" by default
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Row 2 (960004 / 华夏兴华混合H) – scale, position, weight, value, rank changed
$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "10.34"
$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "89.73"
$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "3.37"
$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.3485"
$q4.Range("H2").Value = 10

# Row 3 (006868 / 华夏科技成长股票) – scale, position, weight, value, rank changed
$q4.Range("D3").NumberFormat = "@"
$q4.Range("D3").Value = "5.05"
$q4.Range("E3").NumberFormat = "@"
$q4.Range("E3").Value = "88.93"
$q4.Range("F3").NumberFormat = "@"
$q4.Range("F3").Value = "3.74"
$q4.Range("G3").NumberFormat = "@"
$q4.Range("G3").Value = "0.1889"
$q4.Range("H3").Value = 10

# Row 4 (519908 / 华夏兴华混合A) – position, weight, rank changed
$q4.Range("E4").NumberFormat = "@"
$q4.Range("E4").Value = "89.73"
$q4.Range("F4").NumberFormat = "@"
$q4.Range("F4").Value = "3.37"
$q4.Range("H4").Value = 10

# ------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row for 2022-Q4
#    at the top of the data (row 2), pushing every other quarter
#    down by one row.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$total.Rows(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.54

$total.Range("A2").Font.Bold = $true
$total.Range("A2").HorizontalAlignment = -4108
$total.Range("A2").VerticalAlignment = -4160
$total.Range("A2").Borders.LineStyle = 1
